# Update "想去人数" (wish-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 44
$ws1.Range("F5").Value = 4963
$ws1.Range("F7").Value = 79
$ws1.Range("F8").Value = 285
$ws1.Range("F9").Value = 42

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 44
$ws4.Range("F9").Value = 4963
$ws4.Range("F11").Value = 79
$ws4.Range("F13").Value = 285
$ws4.Range("F14").Value = 42
